$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Point 2.3: accumulating and "averaging" ---

# Row 30: step header values 0,10,20,...,90 across G:P
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 10
$ws.Range("I30").Value = 20
$ws.Range("J30").Value = 30
$ws.Range("K30").Value = 40
$ws.Range("L30").Value = 50
$ws.Range("M30").Value = 60
$ws.Range("N30").Value = 70
$ws.Range("O30").Value = 80
$ws.Range("P30").Value = 90

# Rows 31-36: base values in F, accumulated sums in G:P
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 100
$ws.Range("F33").Value = 200
$ws.Range("F34").Value = 300
$ws.Range("F35").Value = 400
$ws.Range("F36").Value = 500

# G31 alone (not part of a fill-down/across group)
$ws.Range("G31").Formula = '=$F31+G$30'

# H31:P36 filled as one block -> shared formula si=13
$ws.Range("H31:P36").Formula = '=$F31+H$30'

# G32:G36 filled down as one block -> shared formula si=14
$ws.Range("G32:G36").Formula = '=$F32+G$30'

# Row 40: integer "bucket" of the accumulated value from row 31
$ws.Range("G40").Formula = '=QUOTIENT(G31,QUOTIENT(500,9))'
$ws.Range("H40:P40").Formula = '=QUOTIENT(H31,QUOTIENT(500,9))'

# Row 41: bucket of row 32
$ws.Range("G41:P41").Formula = '=QUOTIENT(G32,QUOTIENT(500,9))'

# Row 42: bucket of row 33
$ws.Range("G42:P42").Formula = '=QUOTIENT(G33,QUOTIENT(500,9))'
$ws.Range("U42").Value = 0
$ws.Range("V42").Formula = '=COUNTIF($G$40:$P$45,U42)'

# Row 43: bucket of row 34
$ws.Range("G43:P43").Formula = '=QUOTIENT(G34,QUOTIENT(500,9))'
$ws.Range("U43").Value = 1
$ws.Range("V43:V50").Formula = '=COUNTIF($G$40:$P$45,U43)'

# Row 44: bucket of row 35
$ws.Range("G44:P44").Formula = '=QUOTIENT(G35,QUOTIENT(500,9))'
$ws.Range("U44").Value = 2

# Row 45: bucket of row 36 - only column G ends up populated
$ws.Range("G45:P45").Formula = '=QUOTIENT(G36,QUOTIENT(500,9))'
$ws.Range("H45:P45").ClearContents()
$ws.Range("U45").Value = 3

# Rows 46-50: remaining histogram buckets
$ws.Range("U46").Value = 4
$ws.Range("U47").Value = 5
$ws.Range("U48").Value = 6
$ws.Range("U49").Value = 7
$ws.Range("U50").Value = 8

# Selection / view matches the author's final state
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("W57").Select()
